# Weekly data refresh for "Hortaliza, Vega Modelo de Temuco - Ramas de apio":
# a new week's record is inserted at row 13, pushing the existing rows 13-28
# down to 14-29 (dimension grows from A1:R28 to A1:R29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 13, shifting rows 13-28 -> 14-29.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with this week's new record.
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(13, 3).Value = "La Araucanía"
$ws.Cells.Item(13, 4).Value = 45159
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = 100112017
$ws.Cells.Item(13, 7).Value = "Ramas de apio"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 75
$ws.Cells.Item(13, 11).Value = 5000
$ws.Cells.Item(13, 12).Value = 5000
$ws.Cells.Item(13, 13).Value = 5000
$ws.Cells.Item(13, 14).Value = "$/paquete"
$ws.Cells.Item(13, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(13, 16).Value = 5000
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Hortaliza"
